# Append two new daily rows (2020-09-26 and 2020-09-27) to the
# "out_vars" historical log sheet, matching the upstream SSA data
# refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows go right after the current last row of data (row 118).
$startRow = 118 + 1

$newRows = @(
    @{ Fecha = "2020-09-26"; Confirmados = 726431; Negativos = 851864; Sospechosos = 86762; Defunciones = 76243; Pct = 24.09 },
    @{ Fecha = "2020-09-27"; Confirmados = 730317; Negativos = 857958; Sospechosos = 82914; Defunciones = 76430; Pct = 24.03 }
)

$r = $startRow
foreach ($row in $newRows) {
    # Column A: plain date-label text (matches the existing "Fecha" column,
    # which stores dates as shared-string text, not real Excel date serials).
    $cellA = $ws.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.Fecha
    # Drop back to the default style now that the text has been accepted,
    # so the cell ends up with no explicit style (like the other rows).
    $cellA.Style = "Normal"

    $ws.Range("B$r").Value = $row.Confirmados
    $ws.Range("C$r").Value = $row.Negativos
    $ws.Range("D$r").Value = $row.Sospechosos
    $ws.Range("E$r").Value = $row.Defunciones
    $ws.Range("F$r").Value = $row.Pct

    # These two rows carry a slightly larger, wrap-text style on the
    # numeric columns, same as the freshly appended data from the source
    # refresh.
    $dataRange = $ws.Range("B$r:F$r")
    $dataRange.Font.Size = 12
    $dataRange.WrapText = $true
    $dataRange.EntireRow.RowHeight = 16

    $r = $r + 1
}

$lastRow = $r - 1

# Keep the view roughly where the new data is, mirroring the author's
# scroll position when they made the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
$ws.Range("I107").Select()
